$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Update cell values to reflect the new "Bible API" data
# (order matters for shared-string table layout: aliança, then avraham, then job)
$ws.Range("C2").Value = "aliança"
$ws.Range("B2").Value = "avraham"
$ws.Range("B3").Value = "job"
$ws.Range("C3").Value = "paciencia"
$ws.Range("B4").Value = "moshe"
$ws.Range("C4").Value = "perseverança"

# Update the selected cell/range shown in the sheet view
$ws.Activate()
$ws.Range("C7").Select()
